# Auto-generated edit script: updates cryptos list price/volume cells
# to match the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.831.32"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.354.93"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.47"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.06"
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.13"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.708.10"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.368.38"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.816.47"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.88"
$ws.Range("E19").Value = "  +8.94%  "
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.92"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.63"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.52"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -6.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.85"
$ws.Range("E25").Value = "  +7.20%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.53"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.72"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  +5.62%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.62"
$ws.Range("E33").Value = "  -8.67%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0899"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  +6.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.57"
$ws.Range("E37").Value = "  -7.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  -6.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.68"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.06"
$ws.Range("E44").Value = "  +7.29%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "93.58"
$ws.Range("E46").Value = "  +22.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.92"
$ws.Range("E47").Value = "  -5.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.59"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  -0.24%  "
